$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 30 (Create Payroll) - add dates
$ws.Range("B29:C29").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)
$ws.Range("B30").Value = 45289
$ws.Range("C30").Value = 45289

# Row 31 (Show payroll history) - add dates
$ws.Range("B29:C29").Copy()
$ws.Range("B31:C31").PasteSpecial(-4122)
$ws.Range("B31").Value = 45290
$ws.Range("C31").Value = 45290

# Row 32 (No payroll exists message) - add dates
$ws.Range("B29:C29").Copy()
$ws.Range("B32:C32").PasteSpecial(-4122)
$ws.Range("B32").Value = 45290
$ws.Range("C32").Value = 45290

# Row 33 (Invalid employee id message) - add dates
$ws.Range("B29:C29").Copy()
$ws.Range("B33:C33").PasteSpecial(-4122)
$ws.Range("B33").Value = 45290
$ws.Range("C33").Value = 45290

# Row 34 (Save success message) - add dates
$ws.Range("B29:C29").Copy()
$ws.Range("B34:C34").PasteSpecial(-4122)
$ws.Range("B34").Value = 45290
$ws.Range("C34").Value = 45290

# New row 35 - relisting after save
$ws.Range("A35").Value = "relisting after save "
$ws.Range("B29:C29").Copy()
$ws.Range("B35:C35").PasteSpecial(-4122)
$ws.Range("B35").Value = 45290
$ws.Range("C35").Value = 45290

# New row 36 - implement clear
$ws.Range("A36").Value = "implement clear"
$ws.Range("B29:C29").Copy()
$ws.Range("B36:C36").PasteSpecial(-4122)
$ws.Range("B36").Value = 45290
$ws.Range("C36").Value = 45290

$excel.CutCopyMode = 0

# Update selection to match the new active cell/selection state
$ws.Range("B36:C36").Select()
